$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force text format so numeric-looking strings (e.g. "0.693", "2.42")
# are not reinterpreted as numbers by Excel when assigned via .Value
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "35.145.32"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.905.19"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  -0.49%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "253.33"
$ws.Range("E5").Value = "  +3.29%  "

$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "0.693"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "1.01"
$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "41.45"
$ws.Range("E8").Value = "  +1.94%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.358"
$ws.Range("E9").Value = "  +3.01%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "52.69"
$ws.Range("E10").Value = "  -0.59%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.0751"
$ws.Range("E11").Value = "  +4.42%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.0983"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "13.08"
$ws.Range("E13").Value = "  +4.64%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.177.26"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.734"
$ws.Range("E15").Value = "  +4.33%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "4.98"
$ws.Range("E16").Value = "  +3.63%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.902.07"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "35.148.29"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "73.52"
$ws.Range("E19").Value = "  +1.94%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0833"
$ws.Range("E20").Value = "  +2.09%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "243.62"
$ws.Range("E21").Value = "  +1.27%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "12.95"
$ws.Range("E22").Value = "  +2.96%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.04"
$ws.Range("E23").Value = "  +5.47%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.01"
$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  +4.93%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  +0.84%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "166.73"
$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "8.56"
$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "18.55"
$ws.Range("E29").Value = "  +1.50%  "

$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "0.129"
$ws.Range("E30").Value = "  -1.12%  "

$ws.Range("B31").Value = "EURNeutrino"
$ws.Range("C31").Value = "https://coinranking.com/coin/7YKHKSdb+eurneutrino-eurn"
$ws.Range("D31").Value = "4.128.45"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "2.06"
$ws.Range("E32").Value = "  +13.44%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.0606"
$ws.Range("E33").Value = "  +7.07%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "4.31"
$ws.Range("E34").Value = "  +4.08%  "

$ws.Range("B35").Value = "TrustWalletToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  +1.74%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "4.20"
$ws.Range("E36").Value = "  +2.57%  "

$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "1.01"
$ws.Range("E37").Value = "  -0.50%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.855"
$ws.Range("E38").Value = "  -5.95%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "2.01"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "103.34"
$ws.Range("E40").Value = "  +15.59%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "17.30"
$ws.Range("E41").Value = "  +7.85%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0216"
$ws.Range("E42").Value = "  +3.32%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "1.10"
$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.0648"
$ws.Range("E44").Value = "  -2.63%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  +0.24%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.320.51"
$ws.Range("E46").Value = "  -2.30%  "

$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "2.75"
$ws.Range("E48").Value = "  -1.38%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "6.59"
$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("B50").Value = "Gas"
$ws.Range("C50").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D50").Value = "12.18"
$ws.Range("E50").Value = "  -3.89%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0746"
$ws.Range("E51").Value = "  +5.62%  "

# Restore original (unstyled/"Normal") cell style now that values are set as text
$dataRange.Style = "Normal"
